# CAV Export fix: rows previously tagged with the long-form
# "NO REQUIREMENT FOR CONTINUOUS DIRECTIONAL CONTROL" text in the
# AIRBALANCE RELATIONSHIP column (K) should use the short-form "NR"
# designation instead (matches the value already used elsewhere in the
# workbook, e.g. column F / G "NR" entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2 through 12 (row 1 is the header, row 13 is blank).
$ws.Range("K2:K12").Value = "NR"
